$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.701.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.655.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.87"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.656"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.50%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -3.33%  "

$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("E13").Value = "  -1.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.132.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.594.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.703.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "

$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +10.38%  "

$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("E26").Value = "  -1.51%  "

$ws.Range("E27").Value = "  +2.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "568.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.04%  "

$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("E33").Value = "  +3.07%  "

$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "155.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0608"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("E44").Value = "  -2.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0254"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0242"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.798"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.44%  "
